$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting the
# existing N/O/P ("Late", "Date"/heading, "Outstanding") columns one to
# the right (-> O/P/Q). This mirrors a manual "Insert Column" in Excel.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()

# The freshly inserted column picks up the width of the column that was
# immediately to its left (what Excel does visually on a column insert).
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with the new selection
# left where the user ended up after inserting the column.
$ws.Activate()
$ws.Range("S9").Select()
